# Sample Project / Main.xlsx - "Rules" sheet
#
# The diff changes cell B11 from the shared string "R40" to the (new)
# shared string "1", while leaving the cell's style untouched (s="23").
#
# A plain `$ws.Range("B11").Value = "1"` would store the value as a
# *number* (because "1" looks numeric), which is not what the source
# workbook has (it keeps it a text/shared-string cell, t="s"). To force
# Excel to keep it as text without disturbing B11's existing number
# format/style, we stage the text value in a scratch cell that has been
# explicitly formatted as Text ("@"), copy it, and paste-special just the
# *value* into B11 - this swaps in the text value while B11 keeps its
# original formatting. The scratch cell/row is then removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$helper = $ws.Range("B13")
$helper.NumberFormat = "@"
$helper.Value = "1"
$helper.Copy()

$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues: value only, keep destination formatting

$excel.CutCopyMode = $false
$helper.EntireRow.Delete()
